$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12 / 13: rename source-voltage labels and add new measured values ---
$ws.Range("A12").Value = "Quellspannung_oben_Steckbrett"
$ws.Range("A13").Value = "Quellspannung_unten_Steckbrett"

# --- Row 9: fix duplicate label (was same as row 8) and its measured value ---
$ws.Range("A9").Value = "U_2Weg_ohneC_WDSM"

# --- Row 4 ---
$ws.Range("B4").Formula = "=26/10"
$ws.Range("D4").Value = 0.5

# --- Row 5 ---
$ws.Range("B5").Formula = "=2.6/10"
$ws.Range("D5").Value = 0.5

# --- Row 6 ---
$ws.Range("B6").Formula = "=9.5"
$ws.Range("D6").Value = 0.5

# --- Row 7 ---
$ws.Range("B7").Value = 10.3
$ws.Range("D7").Value = 0.5

# --- Row 8 ---
$ws.Range("B8").Formula = "=61.5/10"
$ws.Range("D8").Value = 0.5

# --- Row 9 (continued) ---
$ws.Range("B9").Formula = "=63.5/10"
$ws.Range("D9").Value = 0.5

# --- Row 10 ---
$ws.Range("B10").Value = 10
$ws.Range("D10").Value = 0.5

# --- Row 11 ---
$ws.Range("B11").Value = 10.3
$ws.Range("D11").Value = 0.5

# --- Row 12 (continued): new measured source voltage ---
$ws.Range("B12").Formula = "=7.1*30"

# --- Row 13 (continued): new measured source voltage ---
$ws.Range("B13").Formula = "=7.5*30"

# --- Update the active selection to D13 ---
[void]$ws.Range("D13").Select()
